{"js": "// Replace the division-problem answers in the table cells with the new values.\n// Each old value is unique within the document, so an exact, case-sensitive\n// search safely targets the single run that needs to change.\nconst replacements = [\n  { oldText: \"639\u00f72=319, 1\", newText: \"421\u00f77=60, 1\" },\n  { oldText: \"108\u00f72=54, 0\", newText: \"536\u00f77=76, 4\" },\n  { oldText: \"378\u00f79=42, 0\", newText: \"699\u00f75=139, 4\" },\n  { oldText: \"519\u00f73=173, 0\", newText: \"412\u00f74=103, 0\" },\n  { oldText: \"818\u00f78=102, 2\", newText: \"115\u00f76=19, 1\" },\n  { oldText: \"206\u00f74=51, 2\", newText: \"465\u00f79=51, 6\" },\n  { oldText: \"235\u00f72=117, 1\", newText: \"999\u00f78=124, 7\" },\n  { oldText: \"271\u00f75=54, 1\", newText: \"648\u00f77=92, 4\" },\n  { oldText: \"330\u00f74=82, 2\", newText: \"367\u00f77=52, 3\" },\n  { oldText: \"798\u00f77=114, 0\", newText: \"669\u00f73=223, 0\" },\n  { oldText: \"556\u00f74=139, 0\", newText: \"128\u00f74=32, 0\" },\n  { oldText: \"572\u00f78=71, 4\", newText: \"714\u00f79=79, 3\" },\n  { oldText: \"200\u00f73=66, 2\", newText: \"977\u00f72=488, 1\" },\n  { oldText: \"981\u00f77=140, 1\", newText: \"919\u00f76=153, 1\" },\n  { oldText: \"315\u00f78=39, 3\", newText: \"980\u00f74=245, 0\" },\n  { oldText: \"808\u00f76=134, 4\", newText: \"644\u00f78=80, 4\" },\n  { oldText: \"980\u00f76=163, 2\", newText: \"319\u00f74=79, 3\" },\n  { oldText: \"857\u00f78=107, 1\", newText: \"367\u00f75=73, 2\" },\n  { oldText: \"551\u00f76=91, 5\", newText: \"729\u00f77=104, 1\" },\n  { oldText: \"489\u00f78=61, 1\", newText: \"909\u00f76=151, 3\" },\n  { oldText: \"902\u00f77=128, 6\", newText: \"894\u00f77=127, 5\" },\n  { oldText: \"534\u00f74=133, 2\", newText: \"182\u00f76=30, 2\" },\n  { oldText: \"186\u00f72=93, 0\", newText: \"392\u00f74=98, 0\" },\n  { oldText: \"393\u00f76=65, 3\", newText: \"827\u00f78=103, 3\" },\n  { oldText: \"560\u00f77=80, 0\", newText: \"655\u00f73=218, 1\" }\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the division-problem answers in the table cells with the new values.\n# Each old value is unique within the document, so Find/Replace with exact,\n# case-sensitive matching safely targets only the intended cell each time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{ Old = \"639\u00f72=319, 1\"; New = \"421\u00f77=60, 1\" },\n  @{ Old = \"108\u00f72=54, 0\"; New = \"536\u00f77=76, 4\" },\n  @{ Old = \"378\u00f79=42, 0\"; New = \"699\u00f75=139, 4\" },\n  @{ Old = \"519\u00f73=173, 0\"; New = \"412\u00f74=103, 0\" },\n  @{ Old = \"818\u00f78=102, 2\"; New = \"115\u00f76=19, 1\" },\n  @{ Old = \"206\u00f74=51, 2\"; New = \"465\u00f79=51, 6\" },\n  @{ Old = \"235\u00f72=117, 1\"; New = \"999\u00f78=124, 7\" },\n  @{ Old = \"271\u00f75=54, 1\"; New = \"648\u00f77=92, 4\" },\n  @{ Old = \"330\u00f74=82, 2\"; New = \"367\u00f77=52, 3\" },\n  @{ Old = \"798\u00f77=114, 0\"; New = \"669\u00f73=223, 0\" },\n  @{ Old = \"556\u00f74=139, 0\"; New = \"128\u00f74=32, 0\" },\n  @{ Old = \"572\u00f78=71, 4\"; New = \"714\u00f79=79, 3\" },\n  @{ Old = \"200\u00f73=66, 2\"; New = \"977\u00f72=488, 1\" },\n  @{ Old = \"981\u00f77=140, 1\"; New = \"919\u00f76=153, 1\" },\n  @{ Old = \"315\u00f78=39, 3\"; New = \"980\u00f74=245, 0\" },\n  @{ Old = \"808\u00f76=134, 4\"; New = \"644\u00f78=80, 4\" },\n  @{ Old = \"980\u00f76=163, 2\"; New = \"319\u00f74=79, 3\" },\n  @{ Old = \"857\u00f78=107, 1\"; New = \"367\u00f75=73, 2\" },\n  @{ Old = \"551\u00f76=91, 5\"; New = \"729\u00f77=104, 1\" },\n  @{ Old = \"489\u00f78=61, 1\"; New = \"909\u00f76=151, 3\" },\n  @{ Old = \"902\u00f77=128, 6\"; New = \"894\u00f77=127, 5\" },\n  @{ Old = \"534\u00f74=133, 2\"; New = \"182\u00f76=30, 2\" },\n  @{ Old = \"186\u00f72=93, 0\"; New = \"392\u00f74=98, 0\" },\n  @{ Old = \"393\u00f76=65, 3\"; New = \"827\u00f78=103, 3\" },\n  @{ Old = \"560\u00f77=80, 0\"; New = \"655\u00f73=218, 1\" }\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.Text = $pair.New\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n  $ok = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)  # wdReplaceAll\n  if (-not $ok) {\n    throw \"Find/Replace failed for: $($pair.Old)\"\n  }\n}\n"}
